# Auto-generated Excel COM-interop script
# Applies numeric updates to the Brynhildr_Profits workbook sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12: Don't Be So Tallow / Beeswax
$ws.Range("H12").Value = 143.5
$ws.Range("I12").Value = 141.91667
$ws.Range("J12").Value = 153
$ws.Range("K12").Value = 141.91667
$ws.Range("L12").Value = 153
$ws.Range("M12").Value = 28.08332999999999
$ws.Range("N12").Value = -493

# Row 16: Using Your Arcane Powers for Fun and Profit / Ash Picatrix
$ws.Range("H16").Value = 209
$ws.Range("I16").Value = 209
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 209
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = 21

# Row 17: One for the Road / Potion
$ws.Range("H17").Value = 1143
$ws.Range("I17").Value = 1056.7
$ws.Range("J17").Value = 1178.9584
$ws.Range("K17").Value = 3170.1
$ws.Range("L17").Value = 3536.8752
$ws.Range("M17").Value = -3002.1
$ws.Range("N17").Value = -3872.8752

# Row 129: Practical Command / Commanding Craftsman's Draught
$ws.Range("H129").Value = 2112.2856
$ws.Range("I129").Value = 1757.6
$ws.Range("J129").Value = 2999
$ws.Range("K129").Value = 5272.799999999999
$ws.Range("L129").Value = 8997
$ws.Range("M129").Value = -272.7999999999993
$ws.Range("N129").Value = -18997

$ws = $wb.Worksheets.Item("ARM")
# Row 4: Eyes Bigger than the Plate / Bronze Plate
$ws.Range("H4").Value = 283.33334
$ws.Range("I4").Value = 283.33334
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 283.33334
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -167.33334

# Row 46: Get Me the Usual / Heavy Steel Flanchard
$ws.Range("H46").Value = 6411.1113
$ws.Range("I46").Value = 3724.75
$ws.Range("J46").Value = 8560.200000000001
$ws.Range("K46").Value = 3724.75
$ws.Range("L46").Value = 8560.200000000001
$ws.Range("M46").Value = -3405.75
$ws.Range("N46").Value = -9198.200000000001

# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 6651.844
$ws.Range("I74").Value = 3593.2432
$ws.Range("J74").Value = 20797.875
$ws.Range("K74").Value = 3593.2432
$ws.Range("L74").Value = 20797.875
$ws.Range("M74").Value = -2719.2432
$ws.Range("N74").Value = -22545.875

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 6651.844
$ws.Range("I77").Value = 3593.2432
$ws.Range("J77").Value = 20797.875
$ws.Range("K77").Value = 17966.216
$ws.Range("L77").Value = 103989.375
$ws.Range("M77").Value = -13598.216
$ws.Range("N77").Value = -112725.375

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 6968.3057
$ws.Range("I132").Value = 4466.4736
$ws.Range("J132").Value = 9764.471
$ws.Range("K132").Value = 13399.4208
$ws.Range("L132").Value = 29293.413
$ws.Range("M132").Value = -10869.4208
$ws.Range("N132").Value = -34353.413

$ws = $wb.Worksheets.Item("BSM")
# Row 13: As Above, Below / Bronze Pickaxe
$ws.Range("H13").Value = 70490
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 70490
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 70490
$ws.Range("N13").Value = -70826

# Row 22: Riveting Run / Iron Rivets
$ws.Range("H22").Value = 283.33334
$ws.Range("I22").Value = 283.33334
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 283.33334
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -110.33334

# Row 26: Unseamly Conditions / Iron Pickaxe
$ws.Range("H26").Value = 29871.5
$ws.Range("I26").Value = 22845.8
$ws.Range("J26").Value = 65000
$ws.Range("K26").Value = 22845.8
$ws.Range("L26").Value = 65000
$ws.Range("M26").Value = -22553.8
$ws.Range("N26").Value = -65584

# Row 86: Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 4230
$ws.Range("I86").Value = 4225
$ws.Range("J86").Value = 4250
$ws.Range("K86").Value = 4225
$ws.Range("L86").Value = 4250
$ws.Range("M86").Value = -3102
$ws.Range("N86").Value = -6496

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 4230
$ws.Range("I89").Value = 4225
$ws.Range("J89").Value = 4250
$ws.Range("K89").Value = 21125
$ws.Range("L89").Value = 21250
$ws.Range("M89").Value = -15509
$ws.Range("N89").Value = -32482

# Row 94: High Steal / High Steel Nugget
$ws.Range("H94").Value = 4149.3887
$ws.Range("I94").Value = 4237.2607
$ws.Range("J94").Value = 3993.923
$ws.Range("K94").Value = 4237.2607
$ws.Range("L94").Value = 3993.923
$ws.Range("M94").Value = -3786.2607
$ws.Range("N94").Value = -4895.923

$ws = $wb.Worksheets.Item("CRP")
# Row 22: Driving Up the Wall / Elm Lumber
$ws.Range("H22").Value = 2006.125
$ws.Range("I22").Value = 540.75
$ws.Range("J22").Value = 3471.5
$ws.Range("K22").Value = 540.75
$ws.Range("L22").Value = 3471.5
$ws.Range("M22").Value = -190.75
$ws.Range("N22").Value = -4171.5

# Row 38: Knock on Wood / Walnut Macuahuitl
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("N38").ClearContents()

# Row 46: Flintstone Fight / Walnut Macuahuitl
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").ClearContents()

# Row 122: Timber of Tenkonto / Horse Chestnut Lumber
$ws.Range("H122").Value = 7808.35
$ws.Range("I122").Value = 1568.7878
$ws.Range("J122").Value = 37223.43
$ws.Range("K122").Value = 4706.3634
$ws.Range("L122").Value = 111670.29
$ws.Range("M122").Value = -2256.3634
$ws.Range("N122").Value = -116570.29

$ws = $wb.Worksheets.Item("CUL")
# Row 86: Let's Not Get Sappy / Birch Syrup
$ws.Range("H86").Value = 386.2
$ws.Range("I86").Value = 461.42856
$ws.Range("J86").Value = 210.66667
$ws.Range("K86").Value = 1384.28568
$ws.Range("L86").Value = 632.00001
$ws.Range("M86").Value = -198.28568
$ws.Range("N86").Value = -3004.00001

# Row 89: Luxury Spillover (L) / Birch Syrup
$ws.Range("H89").Value = 386.2
$ws.Range("I89").Value = 461.42856
$ws.Range("J89").Value = 210.66667
$ws.Range("K89").Value = 4152.85704
$ws.Range("L89").Value = 1896.00003
$ws.Range("M89").Value = 1775.14296
$ws.Range("N89").Value = -13752.00003

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell / Hardsilver Ingot
$ws.Range("H80").Value = 3327
$ws.Range("I80").Value = 2904.5
$ws.Range("J80").Value = 3749.5
$ws.Range("K80").Value = 2904.5
$ws.Range("L80").Value = 3749.5
$ws.Range("M80").Value = -1906.5
$ws.Range("N80").Value = -5745.5

# Row 83: With a Noise That Reaches Heaven (L) / Hardsilver Ingot
$ws.Range("H83").Value = 3327
$ws.Range("I83").Value = 2904.5
$ws.Range("J83").Value = 3749.5
$ws.Range("K83").Value = 14522.5
$ws.Range("L83").Value = 18747.5
$ws.Range("M83").Value = -9530.5
$ws.Range("N83").Value = -28731.5

# Row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 3138.923
$ws.Range("I122").Value = 2942.7368
$ws.Range("J122").Value = 3671.4285
$ws.Range("K122").Value = 8828.2104
$ws.Range("L122").Value = 11014.2855
$ws.Range("M122").Value = -6378.2104
$ws.Range("N122").Value = -15914.2855

# Row 126: Gold Rush Order / Phrygian Gold Ingot
$ws.Range("H126").Value = 2653.077
$ws.Range("I126").Value = 2256
$ws.Range("J126").Value = 3116.3333
$ws.Range("K126").Value = 6768
$ws.Range("L126").Value = 9348.999899999999
$ws.Range("M126").Value = -4298
$ws.Range("N126").Value = -14288.9999

# Row 127: Sage with the Golden Earrings / Phrygian Ear Cuffs of Healing
$ws.Range("H127").Value = 93413
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 93413
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 93413
$ws.Range("N127").Value = -103333

# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 19614.834
$ws.Range("I132").Value = 19614.834
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 58844.50199999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -56314.50199999999

# Row 136: Shiny and Good / Pink Beryl
$ws.Range("H136").Value = 54329.89
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 54329.89
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 162989.67
$ws.Range("N136").Value = -168089.67

$ws = $wb.Worksheets.Item("LTW")
# Row 40: Best Served Toad / Toad Leather
$ws.Range("H40").Value = 4299.7334
$ws.Range("I40").Value = 3735.923
$ws.Range("J40").Value = 7964.5
$ws.Range("K40").Value = 3735.923
$ws.Range("L40").Value = 7964.5
$ws.Range("M40").Value = -3599.923
$ws.Range("N40").Value = -8236.5

# Row 61: Spelling Me Softly / Raptor Leather
$ws.Range("H61").Value = 6269.073
$ws.Range("I61").Value = 6015.7144
$ws.Range("J61").Value = 7747
$ws.Range("K61").Value = 6015.7144
$ws.Range("L61").Value = 7747
$ws.Range("M61").Value = -5813.7144
$ws.Range("N61").Value = -8151

# Row 82: Trainin' the Neck / Dragon Leather
$ws.Range("H82").Value = 3682.9092
$ws.Range("I82").Value = 3682.9092
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 3682.9092
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -3321.9092

# Row 85: Training Is Only Skintight (L) / Dragon Leather
$ws.Range("H85").Value = 3682.9092
$ws.Range("I85").Value = 3682.9092
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 3682.9092
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -2434.9092

# Row 113: Peace in Rest / Atrociraptor Leather
$ws.Range("H113").Value = 6269.073
$ws.Range("I113").Value = 6015.7144
$ws.Range("J113").Value = 7747
$ws.Range("K113").Value = 6015.7144
$ws.Range("L113").Value = 7747
$ws.Range("M113").Value = -3845.7144
$ws.Range("N113").Value = -12087

$ws = $wb.Worksheets.Item("WVR")
# Row 16: Keep It under Wraps / Cotton Turban
$ws.Range("H16").Value = 68255
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 68255
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 68255
$ws.Range("N16").Value = -68839

# Row 107: Flax Wax / Bright Linen Yarn
$ws.Range("H107").Value = 1760.1389
$ws.Range("I107").Value = 1072.0416
$ws.Range("J107").Value = 3136.3333
$ws.Range("K107").Value = 3216.1248
$ws.Range("L107").Value = 9408.999899999999
$ws.Range("M107").Value = -1296.1248
$ws.Range("N107").Value = -13248.9999

# Row 122: Heavy Armoire / Dark Hempen Cloth
$ws.Range("H122").Value = 38822.84
$ws.Range("I122").Value = 2595.4348
$ws.Range("J122").Value = 142976.62
$ws.Range("K122").Value = 7786.3044
$ws.Range("L122").Value = 428929.86
$ws.Range("M122").Value = -5336.3044
$ws.Range("N122").Value = -433829.86

# Row 125: Color Coated / Almasty Serge Coat of Healing
$ws.Range("H125").Value = 58570.57
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 58570.57
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 58570.57
$ws.Range("N125").Value = -68410.57000000001

# Row 126: A Polished Purchase / Snow Linen
$ws.Range("H126").Value = 1482.0454
$ws.Range("I126").Value = 1287.0667
$ws.Range("J126").Value = 1899.8572
$ws.Range("K126").Value = 3861.2001
$ws.Range("L126").Value = 5699.571599999999
$ws.Range("M126").Value = -1391.2001
$ws.Range("N126").Value = -10639.5716
